$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "69.05")
# must have NumberFormat forced to text ("@") before assignment, otherwise
# Excel auto-converts them to numbers and the "xx.xx" 2-decimal text look is lost.
$ws.Range("D2").Value = "65.725.90"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.171.16"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.44"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.34"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").Value = "3.170.54"
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "3.717.63"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.57"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "65.722.41"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "3.175.27"
$ws.Range("E18").Value = "  -4.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.53"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.05"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.495"
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("D26").Value = "3.307.61"
$ws.Range("E27").Value = "  -6.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.84"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.42"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.10"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.51"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.37"
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("D43").Value = "2.656.97"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.17"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.70"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "332.11"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0658"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.00"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0276"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -1.20%  "
